# Generate Report for Handoff
# Adds a new localized file ("ffff27369140-abc4-4010-88af-940ba64a7fee.md")
# to the handoff report, and refreshes the handoff details (uuid/hash/
# timestamps) for the already-tracked file
# ("66098900-5720-41b3-9af7-5ede14027edf.md" -> "465496d9-901d-4e13-a66e-e96712b17117.md").

$wb = $excel.ActiveWorkbook

$oldUuid = "66098900-5720-41b3-9af7-5ede14027edf"
$newUuid = "465496d9-901d-4e13-a66e-e96712b17117"
$newUuid2 = "ffff27369140-abc4-4010-88af-940ba64a7fee"
$newHash = "5686d97c36b38c165bf6378a8fbc03f2a531c9f4"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/f4e838aae3bee7b3f7cd4a9c6958a30ceecc4bd5/e2e"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/f4e838aae3bee7b3f7cd4a9c6958a30ceecc4bd5/.localization-config"
$zhXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/435582a7060e3f19c4fe6d93ca02181c4b0eb248/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
$deXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3353f81e748c0b0d94a8ea99f892cba6255be5d6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"

$newMdName = "$newUuid.md"
$newMdName2 = "$newUuid2.md"
$cfgName = ".localization-config"

$zhXlfName = "$newUuid.$newHash.zh-cn.xlf"
$deXlfName = "$newUuid.$newHash.de-de.xlf"

$readyStatus = "Ready for handoff"
$notLocalized = "Not to be localized"
$includeStatus = "Include"
$ignoredStatus = "Ignored"
$epoch = "0001-01-01 00:00:00"
$zhHandoffTime = "2016-03-04 06:24:44"
$deHandoffTime = "2016-03-04 06:24:59"
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# The bulk "clear all hyperlinks on sheet" quirk; re-add them all below in
# their final, correct order/content.
$ws1.Range("A1").Hyperlinks.Delete()

$ws1.Range("B2").Value = $readyStatus
$ws1.Range("C2").Value = $readyStatus
$ws1.Hyperlinks.Add($ws1.Range("A2"), "$mdBase/$newMdName", "", "", $newMdName) | Out-Null

$ws1.Range("B3").Value = $readyStatus
$ws1.Range("C3").Value = $readyStatus
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$mdBase/$newMdName2", "", "", $newMdName2) | Out-Null

$ws1.Range("B4").Value = $notLocalized
$ws1.Range("C4").Value = $notLocalized
$ws1.Hyperlinks.Add($ws1.Range("A4"), $cfgUrl, "", "", $cfgName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Hyperlinks.Delete()

$ws2.Range("B2").Value = $readyStatus
$ws2.Range("D2").Value = $zhHandoffTime
$ws2.Range("G2").Value = $epoch
$ws2.Range("H2").Value = $includeStatus
$ws2.Hyperlinks.Add($ws2.Range("A2"), "$mdBase/$newMdName", "", "", $newMdName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$zhXlfBase/$zhXlfName", "", "", $zhXlfName) | Out-Null

$ws2.Range("B3").Value = $readyStatus
$ws2.Range("D3").Value = $zhHandoffTime
$ws2.Range("G3").Value = $epoch
$ws2.Range("H3").Value = $includeStatus
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$mdBase/$newMdName2", "", "", $newMdName2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "$zhXlfBase/$zhXlfName", "", "", $zhXlfName) | Out-Null

$ws2.Range("B4").Value = $notLocalized
$ws2.Range("D4").Value = $epoch
$ws2.Range("D4").NumberFormat = $dateFmt
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = $ignoredStatus
$ws2.Hyperlinks.Add($ws2.Range("A4"), $cfgUrl, "", "", $cfgName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A1").Hyperlinks.Delete()

$ws3.Range("B2").Value = $readyStatus
$ws3.Range("D2").Value = $deHandoffTime
$ws3.Range("G2").Value = $epoch
$ws3.Range("H2").Value = $includeStatus
$ws3.Hyperlinks.Add($ws3.Range("A2"), "$mdBase/$newMdName", "", "", $newMdName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$deXlfBase/$deXlfName", "", "", $deXlfName) | Out-Null

$ws3.Range("B3").Value = $readyStatus
$ws3.Range("D3").Value = $deHandoffTime
$ws3.Range("G3").Value = $epoch
$ws3.Range("H3").Value = $includeStatus
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$mdBase/$newMdName2", "", "", $newMdName2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "$deXlfBase/$deXlfName", "", "", $deXlfName) | Out-Null

$ws3.Range("B4").Value = $notLocalized
$ws3.Range("D4").Value = $epoch
$ws3.Range("D4").NumberFormat = $dateFmt
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = $ignoredStatus
$ws3.Hyperlinks.Add($ws3.Range("A4"), $cfgUrl, "", "", $cfgName) | Out-Null

Write-Host "Handoff report regenerated: added $newMdName2 and refreshed $newMdName across Overview/zh-cn/de-de sheets."
